# Build ANOVA_df for TRA: remove all "除夕" (New Year's Eve) holiday rows
# from the Holidays table, since that category is not used in the ANOVA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (as currently numbered, top to bottom) whose category column (C)
# is the "除夕" shared string. Deleting from the bottom up keeps the
# remaining row numbers valid while we iterate.
$rowsToDelete = @(10, 19, 28, 37, 46, 55, 64, 73, 82, 91, 100, 109, 118, 127, 136, 145, 154, 163)

for ($i = $rowsToDelete.Count - 1; $i -ge 0; $i--) {
    $r = $rowsToDelete[$i]
    $ws.Rows($r).EntireRow.Delete()
}

# Restore the view/selection state recorded after the edit.
$ws.Range("F139").Select()
try {
    $excel.ActiveWindow.ScrollRow = 121
} catch {
}

Write-Host "Removed $($rowsToDelete.Count) holiday rows from Sheet1."
